$d = $word.ActiveDocument

$found = $d.Content.Find.Execute(
    "NTA AP Invoice Confirm screen is used to complete the NTA AP Invoice. The NTA AP Invoice transaction inprogress will be shown here, User have to complete the NTA  AP Invoice Confirm inorder to complete the NTA AP Invoice.",
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    $rng = $d.Content.Duplicate
    $rng.Start = $d.Content.Find.Parent.Start
}

Write-Output "Found: $found"
